# Weekly update for the Albahaca (Mercado Mayorista Lo Valledor de Santiago) sheet.
# Two brand-new daily records are inserted at the top of the data block (rows 233-234),
# pushing the previously-recorded rows down by two; the two oldest rows that fall off
# the bottom of the original block (old rows 261-262) are appended as new rows 263-264.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 233
$lastRow  = 262

# --- 1. Snapshot the "before" state of the variable columns (D, I, J-Q) for every
#        row in the block, before we overwrite anything. ---
$snapD = @{}
$snapI = @{}
$snapJ = @{}
$snapK = @{}
$snapL = @{}
$snapM = @{}
$snapN = @{}
$snapO = @{}
$snapP = @{}
$snapQ = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapD[$r] = $ws.Cells.Item($r, 4).Value2
    $snapI[$r] = $ws.Cells.Item($r, 9).Value2
    $snapJ[$r] = $ws.Cells.Item($r, 10).Value2
    $snapK[$r] = $ws.Cells.Item($r, 11).Value2
    $snapL[$r] = $ws.Cells.Item($r, 12).Value2
    $snapM[$r] = $ws.Cells.Item($r, 13).Value2
    $snapN[$r] = $ws.Cells.Item($r, 14).Value2
    $snapO[$r] = $ws.Cells.Item($r, 15).Value2
    $snapP[$r] = $ws.Cells.Item($r, 16).Value2
    $snapQ[$r] = $ws.Cells.Item($r, 17).Value2
}

# --- 2. Shift rows 235..262 down from what used to be rows 233..260. ---
for ($r = $lastRow; $r -ge ($firstRow + 2); $r--) {
    $src = $r - 2
    $ws.Cells.Item($r, 4).Value2  = $snapD[$src]
    $ws.Cells.Item($r, 9).Value2  = $snapI[$src]
    $ws.Cells.Item($r, 10).Value2 = $snapJ[$src]
    $ws.Cells.Item($r, 11).Value2 = $snapK[$src]
    $ws.Cells.Item($r, 12).Value2 = $snapL[$src]
    $ws.Cells.Item($r, 13).Value2 = $snapM[$src]
    $ws.Cells.Item($r, 14).Value2 = $snapN[$src]
    $ws.Cells.Item($r, 15).Value2 = $snapO[$src]
    $ws.Cells.Item($r, 16).Value2 = $snapP[$src]
    $ws.Cells.Item($r, 17).Value2 = $snapQ[$src]
}

# --- 3. Write the two brand-new records into rows 233-234. ---
$ws.Cells.Item(233, 4).Value2  = 44505
$ws.Cells.Item(233, 9).Value2  = "Primera"
$ws.Cells.Item(233, 10).Value2 = 350
$ws.Cells.Item(233, 11).Value2 = 5500
$ws.Cells.Item(233, 12).Value2 = 6000
$ws.Cells.Item(233, 13).Value2 = 5786
$ws.Cells.Item(233, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(233, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(233, 16).Value2 = 964
$ws.Cells.Item(233, 17).Value2 = 6

$ws.Cells.Item(234, 4).Value2  = 44505
$ws.Cells.Item(234, 9).Value2  = "Segunda"
$ws.Cells.Item(234, 10).Value2 = 100
$ws.Cells.Item(234, 11).Value2 = 5000
$ws.Cells.Item(234, 12).Value2 = 5000
$ws.Cells.Item(234, 13).Value2 = 5000
$ws.Cells.Item(234, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(234, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(234, 16).Value2 = 833
$ws.Cells.Item(234, 17).Value2 = 6

# --- 4. Append the two rows that fell off the bottom (old rows 261-262) as new
#        rows 263-264, copying every column (A-R) since these rows did not exist
#        before. ---
$newRows = @(263, 264)
$oldRows = @(261, 262)

for ($i = 0; $i -lt 2; $i++) {
    $nr = $newRows[$i]
    $or = $oldRows[$i]

    $ws.Cells.Item($nr, 1).Value2  = $ws.Cells.Item($or, 1).Value2   # Mercado ID
    $ws.Cells.Item($nr, 2).Value2  = $ws.Cells.Item($or, 2).Value2   # Mercado
    $ws.Cells.Item($nr, 3).Value2  = $ws.Cells.Item($or, 3).Value2   # Region
    $ws.Cells.Item($nr, 4).Value2  = $snapD[$or]                     # Fecha
    $ws.Cells.Item($nr, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($nr, 5).Value2  = $ws.Cells.Item($or, 5).Value2   # Codreg
    $ws.Cells.Item($nr, 6).Value2  = $ws.Cells.Item($or, 6).Value2   # Categoria ID
    $ws.Cells.Item($nr, 7).Value2  = $ws.Cells.Item($or, 7).Value2   # Categoria
    $ws.Cells.Item($nr, 8).Value2  = $ws.Cells.Item($or, 8).Value2   # Variedad
    $ws.Cells.Item($nr, 9).Value2  = $snapI[$or]                     # Calidad
    $ws.Cells.Item($nr, 10).Value2 = $snapJ[$or]                     # Volumen
    $ws.Cells.Item($nr, 11).Value2 = $snapK[$or]                     # Precio minimo
    $ws.Cells.Item($nr, 12).Value2 = $snapL[$or]                     # Precio maximo
    $ws.Cells.Item($nr, 13).Value2 = $snapM[$or]                     # Precio promedio ponderado
    $ws.Cells.Item($nr, 14).Value2 = $snapN[$or]                     # Unidad de comercializacion
    $ws.Cells.Item($nr, 15).Value2 = $snapO[$or]                     # Origen
    $ws.Cells.Item($nr, 16).Value2 = $snapP[$or]                     # Precio $/Kg
    $ws.Cells.Item($nr, 17).Value2 = $snapQ[$or]                     # Kg o Unidades
    $ws.Cells.Item($nr, 18).Value2 = $ws.Cells.Item($or, 18).Value2  # Clasificacion
}
